# Update the concise marksheet: correct/total marks changes on sheet "quiz"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# B11: Marking row - "Right" column value changes from 3 to 5
$ws.Range("B11").Value = 5

# B12: Total row - "Right" column value changes from 63 to 105
$ws.Range("B12").Value = 105

# E12: Total row - "Max" column text changes from "62/84" to "105/140"
$ws.Range("E12").Value = "105/140"
